# brake jerk related kpis extracted
#
# Summary of the edit:
#  - params sheet: row 14 (JERK_THD) becomes JERK_NEG_THD with a value and new
#    description; three new parameter rows are appended (JERK_POS_THD,
#    BRAKEJERK_MIN_SPEED, BRAKEJERK_MAX_SPEED); column A is widened.
#  - KPI sheet: the brakeJerkThd row is removed, and the former brakeAccelMax
#    row is renamed to brakeAccelMin.
#  - The params sheet becomes the active / selected sheet instead of
#    vbRcSignals, and a couple of sheets remember a new selected cell.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "params" sheet - brake-jerk parameter rows
# ---------------------------------------------------------------------
$params = $wb.Worksheets.Item("params")

# Row 14 used to describe a single fixed jerk threshold (JERK_THD); it is
# repurposed as the negative-jerk threshold and gets a real default value.
$params.Range("A14").Value = "JERK_NEG_THD"
$params.Range("B14").Value = -20
$params.Range("C14").Value = "float"
$params.Range("D14").Value = "m/s³"
$params.Range("E14").Value = "Negative jerk threshold"
$params.Range("F14").Value = "FcwKpiExtractor"

# New row: positive-jerk threshold.
$params.Range("A15").Value = "JERK_POS_THD"
$params.Range("B15").Value = 20
$params.Range("C15").Value = "float"
$params.Range("D15").Value = "m/s³"
$params.Range("E15").Value = "Positive jerk threshold"
$params.Range("F15").Value = "FcwKpiExtractor"

# New row: minimum speed gate for the brake-jerk warning.
$params.Range("A16").Value = "BRAKEJERK_MIN_SPEED"
$params.Range("B16").Value = 30
$params.Range("C16").Value = "float"
$params.Range("D16").Value = "kph"
$params.Range("E16").Value = "minimum speed threshold for brake jerk warning"
$params.Range("F16").Value = "FcwKpiExtractor"

# New row: maximum speed gate for the brake-jerk warning.
$params.Range("A17").Value = "BRAKEJERK_MAX_SPEED"
$params.Range("B17").Value = 130
$params.Range("C17").Value = "float"
$params.Range("D17").Value = "kph"
$params.Range("E17").Value = "maximum speed threshold for brake jerk warning"
$params.Range("F17").Value = "FcwKpiExtractor"

# Column A needs to widen to fit the longer parameter names.
$params.Columns.Item(1).ColumnWidth = 21.5

# ---------------------------------------------------------------------
# 2. "KPI" sheet - drop brakeJerkThd, rename brakeAccelMax -> brakeAccelMin
# ---------------------------------------------------------------------
$kpi = $wb.Worksheets.Item("KPI")

# The brakeJerkThd row (row 40) is no longer produced, remove it entirely;
# remaining rows shift up automatically.
$kpi.Rows.Item(40).Delete()

# The former brakeAccelMax row (now row 41 after the delete above) becomes
# brakeAccelMin; everything else about that row stays the same.
$kpi.Range("C41").Value = "brakeAccelMin"

# ---------------------------------------------------------------------
# 3. Selection / active sheet bookkeeping
# ---------------------------------------------------------------------
# KPI sheet remembers a new selected cell in its frozen pane.
$kpi.Activate()
$kpi.Range("I28").Select()

# "params" becomes the active / selected sheet (was "vbRcSignals").
$params.Activate()
$params.Range("E25").Select()
